$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Fix <cols>: column A should only span col 1 (was col 1-2),
#    so column B falls back to its own (wider) width/style rule.
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.8

# ---------------------------------------------------------------
# 2) Insert a new row at 12 for the "Docentes responsaveis:" label
#    (shifts old rows 12-21 down to 13-22).
# ---------------------------------------------------------------
$ws.Rows.Item(12).Insert()

# ---------------------------------------------------------------
# 3) Write the corrected cell text for every row.
#    Numeric-looking literals ("2", "0", "01/01/2021") are pushed
#    through a scratch formula-cell + paste-values so Excel keeps
#    them as plain text instead of coercing them to number/date.
# ---------------------------------------------------------------
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOQ4267"
$ws.Range("C2").Value = "LOQ4267"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Tópicos Especiais em Engenharia de Produção II"
$ws.Range("C3").Value = " Tópicos Especiais em Engenharia de Produção II"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Special Topics in Production Engineering II"
$ws.Range("C4").Value = "Special Topics in Production Engineering II"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("Z1").Formula = "=""2"""
$ws.Range("Z1").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("Z1").Formula = "=""2"""
$ws.Range("Z1").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("Z1").Formula = "=""0"""
$ws.Range("Z1").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("Z1").Formula = "=""0"""
$ws.Range("Z1").Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("Z1").Formula = "=""01/01/2021"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("Z1").Formula = "=""01/01/2021"""
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial(-4163)

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EP-9"
$ws.Range("C9").Value = "EP-9"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("C10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"
$ws.Range("C11").Value = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"
$ws.Rows.Item(11).RowHeight = 60

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "A definir de acordo com o tópico programado"
$ws.Range("C14").Value = "A definir de acordo com o tópico programado"
$ws.Rows.Item(14).RowHeight = 60

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "To be defined according to the scheduled topic"
$ws.Range("C15").Value = "To be defined according to the scheduled topic"
$ws.Rows.Item(15).RowHeight = 60

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("C16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Rows.Item(16).RowHeight = 120

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C17").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Rows.Item(17).RowHeight = 120

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Rows.Item(19).RowHeight = 60

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."
$ws.Range("C22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."
$ws.Rows.Item(22).RowHeight = 120

# Clear the scratch cell used for forcing text literals
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------
# 4) Remove the stray "Docentes responsaveis:" text that the row
#    insert/shift left behind in column A of row 13 (that label
#    now correctly lives in A12 instead).
# ---------------------------------------------------------------
$ws.Range("A13").ClearContents()

# ---------------------------------------------------------------
# 5) Safety pass: re-assert the per-column style (font/alignment)
#    for every cell in the shifted rows 12-22, in case the engine
#    mis-assigned a style when materialising a previously-empty
#    cell. Row 3 holds the canonical style for each column.
# ---------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)

Write-Host "Edit applied successfully"